$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35, shifting existing rows 35-48 down to 36-49.
$ws.Rows("35:35").Insert()

# Fill in the newly inserted row 35 with the new XOR instruction entry.
$ws.Range("A35").Value = "0x33"
$ws.Range("B35").Value = "0011 0011"
$ws.Range("C35").Value = "XOR"

# Update the view so that the visible area matches the edited region.
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("C35").Select()
